$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-10 (columns A-F)
$data = @(
    @(27, 0.8107328646344186, -1.404552951823922, 0.2791825964964812, 0.4642850981975821, 0.3576535517144857),
    @(28, 1.197979265296577, -1.383524806554267, 0.2793392705239028, 0.4611815990091802, 0.3560022500362016),
    @(29, 5.555209259952779, -1.196796697727143, 0.295648926262209, 0.4756433605480084, 0.3792873153245865),
    @(30, 2.024697424499826, -1.349387073041513, 0.283575169927707, 0.457505340064679, 0.3643838030461788),
    @(31, 0.9735009822949687, -1.39625761291196, 0.2752441701796778, 0.4556972671299721, 0.3559362525705976),
    @(32, -0.9126926661758791, -1.48086255930799, 0.2741533070165214, 0.4606079045795512, 0.3501436751722561),
    @(33, 2.231027972550454, -1.341482471523983, 0.2807881996020219, 0.4612579830966981, 0.3624319567207067),
    @(34, 0.2549752556816626, -1.425702310784218, 0.2751145571388615, 0.4613300348097312, 0.3534310911770412),
    @(35, 4.474318729867433, -1.239134840997663, 0.2877082482069251, 0.4638503111729836, 0.3705068175371698)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}

# Apply the same formatting used for column A (rows 2-5) to the new column A cells (rows 6-10)
$ws.Range("A2").Copy()
$ws.Range("A6:A10").PasteSpecial(-4122)
$excel.CutCopyMode = 0
